$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 519.1177
$ws.Range("I92").Value = 405.33334
$ws.Range("J92").Value = 1372.5
$ws.Range("K92").Value = 405.33334
$ws.Range("L92").Value = 1372.5
$ws.Range("M92").Value = 842.66666
$ws.Range("N92").Value = -3868.5

$ws.Range("H98").Value = 2933.0908
$ws.Range("I98").Value = 1472.7142
$ws.Range("J98").Value = 5488.75
$ws.Range("K98").Value = 1472.7142
$ws.Range("L98").Value = 5488.75
$ws.Range("M98").Value = 25.28580000000011
$ws.Range("N98").Value = -8484.75

$ws.Range("H122").Value = 2933.0908
$ws.Range("I122").Value = 1472.7142
$ws.Range("J122").Value = 5488.75
$ws.Range("K122").Value = 4418.142599999999
$ws.Range("L122").Value = 16466.25
$ws.Range("M122").Value = -1968.142599999999
$ws.Range("N122").Value = -21366.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4486.303
$ws.Range("I61").Value = 4459.625
$ws.Range("J61").Value = 4557.4443
$ws.Range("K61").Value = 4459.625
$ws.Range("L61").Value = 4557.4443
$ws.Range("M61").Value = -4247.625
$ws.Range("N61").Value = -4981.4443

$ws.Range("H74").Value = 2816.9688
$ws.Range("I74").Value = 2019
$ws.Range("J74").Value = 4146.9165
$ws.Range("K74").Value = 2019
$ws.Range("L74").Value = 4146.9165
$ws.Range("M74").Value = -1145
$ws.Range("N74").Value = -5894.9165

$ws.Range("H77").Value = 2816.9688
$ws.Range("I77").Value = 2019
$ws.Range("J77").Value = 4146.9165
$ws.Range("K77").Value = 10095
$ws.Range("L77").Value = 20734.5825
$ws.Range("M77").Value = -5727
$ws.Range("N77").Value = -29470.5825

$ws.Range("H97").Value = 1568.8438
$ws.Range("I97").Value = 1139.1666
$ws.Range("K97").Value = 1139.1666
$ws.Range("M97").Value = -643.1666

$ws.Range("H131").Value = 79997.75
$ws.Range("J131").Value = 79997.75
$ws.Range("L131").Value = 79997.75
$ws.Range("N131").Value = -90077.75

$ws.Range("H132").Value = 4106.364
$ws.Range("I132").Value = 4159.048
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 12477.144
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -9947.144
$ws.Range("N132").Value = -14060

$ws.Range("H136").Value = 4486.303
$ws.Range("I136").Value = 4459.625
$ws.Range("J136").Value = 4557.4443
$ws.Range("K136").Value = 13378.875
$ws.Range("L136").Value = 13672.3329
$ws.Range("M136").Value = -10828.875
$ws.Range("N136").Value = -18772.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2648.6
$ws.Range("I107").Value = 2078.818
$ws.Range("K107").Value = 2078.818
$ws.Range("M107").Value = -158.8180000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 144999.5
$ws.Range("I16").Value = 150000
$ws.Range("J16").Value = 139999
$ws.Range("K16").Value = 150000
$ws.Range("L16").Value = 139999
$ws.Range("M16").Value = -149713
$ws.Range("N16").Value = -140573

$ws.Range("H105").Value = 1775.25
$ws.Range("I105").Value = 1700.3334
$ws.Range("K105").Value = 1700.3334
$ws.Range("M105").Value = 46.66660000000002

$ws.Range("H113").Value = 144999.5
$ws.Range("I113").Value = 150000
$ws.Range("J113").Value = 139999
$ws.Range("K113").Value = 150000
$ws.Range("L113").Value = 139999
$ws.Range("M113").Value = -147830
$ws.Range("N113").Value = -144339

$ws.Range("H124").Value = 36800
$ws.Range("J124").Value = 48600
$ws.Range("L124").Value = 48600
$ws.Range("N124").Value = -53510

$ws.Range("H132").Value = 3969.625
$ws.Range("I132").Value = 3129.5715
$ws.Range("K132").Value = 9388.7145
$ws.Range("M132").Value = -6858.7145

$ws.Range("H134").Value = 2984.3333
$ws.Range("I134").Value = 2984.3333
$ws.Range("K134").Value = 8952.999899999999
$ws.Range("M134").Value = -6417.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 192.55
$ws.Range("I107").Value = 172.83333
$ws.Range("J107").Value = 222.125
$ws.Range("K107").Value = 172.83333
$ws.Range("L107").Value = 222.125
$ws.Range("M107").Value = 1747.16667
$ws.Range("N107").Value = -4062.125

$ws.Range("H122").Value = 4234.9443
$ws.Range("I122").Value = 1698.6364
$ws.Range("J122").Value = 8220.571
$ws.Range("K122").Value = 5095.9092
$ws.Range("L122").Value = 24661.713
$ws.Range("M122").Value = -2645.9092
$ws.Range("N122").Value = -29561.713

$ws.Range("H123").Value = 34999.4
$ws.Range("J123").Value = 34999.4
$ws.Range("L123").Value = 34999.4
$ws.Range("N123").Value = -39899.4

$ws.Range("H132").Value = 3107.3333
$ws.Range("I132").Value = 2745.1428
$ws.Range("J132").Value = 4375
$ws.Range("K132").Value = 8235.428400000001
$ws.Range("L132").Value = 13125
$ws.Range("M132").Value = -5705.428400000001
$ws.Range("N132").Value = -18185

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 6070.857
$ws.Range("I25").Value = 5416.1665
$ws.Range("J25").Value = 9999
$ws.Range("K25").Value = 5416.1665
$ws.Range("L25").Value = 9999
$ws.Range("M25").Value = -5186.1665
$ws.Range("N25").Value = -10459

$ws.Range("H40").Value = 23101.6
$ws.Range("I40").Value = 3668
$ws.Range("J40").Value = 52252
$ws.Range("K40").Value = 3668
$ws.Range("L40").Value = 52252
$ws.Range("M40").Value = -3532
$ws.Range("N40").Value = -52524

$ws.Range("H100").Value = 10000
$ws.Range("I100").Value = 10000
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 10000
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -9459
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value = 5597.8
$ws.Range("I122").Value = 5597.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 16793.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -14343.4
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 6889.143
$ws.Range("I136").Value = 6889.5
$ws.Range("K136").Value = 20668.5
$ws.Range("M136").Value = -18118.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3611.8
$ws.Range("I132").Value = 2030
$ws.Range("K132").Value = 6090
$ws.Range("M132").Value = -3560
